$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 18
$ws.Range("H18").Value = 1473.6666
$ws.Range("I18").Value = 1368
$ws.Range("J18").Value = 2002
$ws.Range("K18").Value = 1368
$ws.Range("L18").Value = 2002
$ws.Range("M18").Value = -1084
$ws.Range("N18").Value = -2570

# Row 34
$ws.Range("H34").Value = 6180
$ws.Range("I34").Value = 6180
$ws.Range("K34").Value = 6180
$ws.Range("M34").Value = -5977

# Row 36
$ws.Range("H36").Value = 6180
$ws.Range("I36").Value = 6180
$ws.Range("K36").Value = 6180
$ws.Range("M36").Value = -5465

# Row 69
$ws.Range("H69").Value = 7777
$ws.Range("J69").Value = 7777
$ws.Range("L69").Value = 23331
$ws.Range("N69").Value = -25079

# Row 72
$ws.Range("H72").Value = 7777
$ws.Range("J72").Value = 7777
$ws.Range("L72").Value = 69993
$ws.Range("N72").Value = -78729

# Row 135
$ws.Range("H135").Value = 855.5833
$ws.Range("I135").Value = 698.2
$ws.Range("K135").Value = 6283.8
$ws.Range("M135").Value = -3748.8

# Row 136
$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 2114.7222
$ws.Range("I61").Value = 1396.5
$ws.Range("K61").Value = 1396.5
$ws.Range("M61").Value = -1184.5

# Row 74
$ws.Range("H74").Value = 6268.7
$ws.Range("I74").Value = 5803.8335
$ws.Range("K74").Value = 5803.8335
$ws.Range("M74").Value = -4929.8335

# Row 77
$ws.Range("H77").Value = 6268.7
$ws.Range("I77").Value = 5803.8335
$ws.Range("K77").Value = 29019.1675
$ws.Range("M77").Value = -24651.1675

# Row 122
$ws.Range("H122").Value = 1257
$ws.Range("I122").Value = 834.4545000000001
$ws.Range("K122").Value = 2503.3635
$ws.Range("M122").Value = -53.36350000000039

# Row 136
$ws.Range("H136").Value = 2114.7222
$ws.Range("I136").Value = 1396.5
$ws.Range("K136").Value = 4189.5
$ws.Range("M136").Value = -1639.5

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 3282.3333
$ws.Range("I20").Value = 2939
$ws.Range("J20").Value = 4999
$ws.Range("K20").Value = 2939
$ws.Range("L20").Value = 4999
$ws.Range("M20").Value = -2692
$ws.Range("N20").Value = -5493

# Row 80
$ws.Range("H80").Value = 860.1818
$ws.Range("I80").Value = 429.33334
$ws.Range("K80").Value = 429.33334
$ws.Range("M80").Value = 568.66666

# Row 83
$ws.Range("H83").Value = 860.1818
$ws.Range("I83").Value = 429.33334
$ws.Range("K83").Value = 2146.6667
$ws.Range("M83").Value = 2845.3333

# Row 86
$ws.Range("H86").Value = 5598.5
$ws.Range("I86").Value = 1798.8
$ws.Range("J86").Value = 8312.571
$ws.Range("K86").Value = 1798.8
$ws.Range("L86").Value = 8312.571
$ws.Range("M86").Value = -675.8
$ws.Range("N86").Value = -10558.571

# Row 89
$ws.Range("H89").Value = 5598.5
$ws.Range("I89").Value = 1798.8
$ws.Range("J89").Value = 8312.571
$ws.Range("K89").Value = 8994
$ws.Range("L89").Value = 41562.855
$ws.Range("M89").Value = -3378
$ws.Range("N89").Value = -52794.855

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 2279
$ws.Range("I22").Value = 1223.3334
$ws.Range("J22").Value = 3070.75
$ws.Range("K22").Value = 1223.3334
$ws.Range("L22").Value = 3070.75
$ws.Range("M22").Value = -873.3334
$ws.Range("N22").Value = -3770.75

# Row 31
$ws.Range("H31").Value = 5354.16
$ws.Range("I31").Value = 1921.4
$ws.Range("J31").Value = 6212.35
$ws.Range("K31").Value = 1921.4
$ws.Range("L31").Value = 6212.35
$ws.Range("M31").Value = -1626.4
$ws.Range("N31").Value = -6802.35

# Row 34
$ws.Range("H34").Value = 5354.16
$ws.Range("I34").Value = 1921.4
$ws.Range("J34").Value = 6212.35
$ws.Range("K34").Value = 1921.4
$ws.Range("L34").Value = 6212.35
$ws.Range("M34").Value = -1719.4
$ws.Range("N34").Value = -6616.35

# Row 39
$ws.Range("H39").Value = 30468
$ws.Range("I39").Value = 2001
$ws.Range("K39").Value = 2001
$ws.Range("M39").Value = -1610

# Row 49
$ws.Range("H49").Value = 30468
$ws.Range("I49").Value = 2001
$ws.Range("K49").Value = 2001
$ws.Range("M49").Value = -1819

$ws = $wb.Worksheets.Item("CUL")
# Row 34
$ws.Range("H34").Value = 2529.9
$ws.Range("I34").Value = 184.33333
$ws.Range("J34").Value = 3535.1428
$ws.Range("K34").Value = 552.99999
$ws.Range("L34").Value = 10605.4284
$ws.Range("M34").Value = -468.99999
$ws.Range("N34").Value = -10773.4284

# Row 54
$ws.Range("H54").Value = 2137.5
$ws.Range("J54").Value = 1500
$ws.Range("L54").Value = 4500
$ws.Range("N54").Value = -5618

# Row 109
$ws.Range("H109").Value = 1292.6666
$ws.Range("I109").Value = 495
$ws.Range("J109").Value = 2888
$ws.Range("K109").Value = 1485
$ws.Range("L109").Value = 8664
$ws.Range("M109").Value = -445
$ws.Range("N109").Value = -10744

# Row 134
$ws.Range("H134").Value = 2208.3333
$ws.Range("I134").Value = 2208.3333
$ws.Range("K134").Value = 6624.999899999999
$ws.Range("M134").Value = -1554.999899999999

$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 336.17648
$ws.Range("J2").Value = 604
$ws.Range("L2").Value = 604
$ws.Range("N2").Value = -830

# Row 36
$ws.Range("H36").Value = 2858.1428
$ws.Range("I36").Value = 1379.25
$ws.Range("J36").Value = 4830
$ws.Range("K36").Value = 1379.25
$ws.Range("L36").Value = 4830
$ws.Range("M36").Value = -894.25
$ws.Range("N36").Value = -5800

# Row 132
$ws.Range("H132").Value = 26500
$ws.Range("I132").Value = 26500
$ws.Range("K132").Value = 79500
$ws.Range("M132").Value = -76970

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1157
$ws.Range("I22").Value = 900
$ws.Range("K22").Value = 900
$ws.Range("M22").Value = -605

# Row 27
$ws.Range("H27").Value = 1157
$ws.Range("I27").Value = 900
$ws.Range("K27").Value = 900
$ws.Range("M27").Value = -793

# Row 40
$ws.Range("H40").Value = 4197.6665
$ws.Range("I40").Value = 3111.2856
$ws.Range("K40").Value = 3111.2856
$ws.Range("M40").Value = -2975.2856

# Row 46
$ws.Range("H46").Value = 6857.7896
$ws.Range("I46").Value = 5024.8335
$ws.Range("K46").Value = 5024.8335
$ws.Range("M46").Value = -4836.8335

# Row 100
$ws.Range("H100").Value = 6982.6
$ws.Range("I100").Value = 3965.2
$ws.Range("J100").Value = 10000
$ws.Range("K100").Value = 3965.2
$ws.Range("L100").Value = 10000
$ws.Range("M100").Value = -3424.2
$ws.Range("N100").Value = -11082

$ws = $wb.Worksheets.Item("WVR")
# Row 96
$ws.Range("H96").Value = 2750
$ws.Range("J96").Value = 4000
$ws.Range("L96").Value = 4000
$ws.Range("N96").Value = -6746
